$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-30 Sunday", "2025-03-31 Monday"),
    @("28÷9=", "98÷6="),
    @("81÷8=", "34÷9="),
    @("31÷5=", "72÷8="),
    @("74÷7=", "10÷4="),
    @("48÷5=", "73÷5="),
    @("96÷3=", "64÷3="),
    @("68÷2=", "48÷9="),
    @("58÷4=", "93÷6="),
    @("41÷4=", "44÷5="),
    @("23÷5=", "81÷2="),
    @("55÷4=", "16÷3="),
    @("50÷3=", "78÷8="),
    @("51÷7=", "98÷8="),
    @("71÷5=", "33÷6="),
    @("43÷4=", "33÷8="),
    @("48÷2=", "75÷3="),
    @("23÷9=", "92÷9="),
    @("46÷4=", "23÷4="),
    @("79÷7=", "97÷3="),
    @("65÷4=", "24÷9="),
    @("23÷8=", "17÷3="),
    @("16÷4=", "92÷8="),
    @("38÷5=", "87÷2="),
    @("77÷6=", "53÷3="),
    @("52÷6=", "10÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
